$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.225.15"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.87%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.915.05"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9977"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "334.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9983"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.24%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4652"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.16%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.72"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08036"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.013"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.27%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.99"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.898.17"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.34%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.971"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.21%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.119"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "89.51"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9978"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.41%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001036"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06585"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.68%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9989"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.19%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "29.240.89"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +2.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.462"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.36"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.62%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.239"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.66%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.121.29"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "157.42"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.80"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.124"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.07%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.452"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.08%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.67"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.42%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9926"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09440"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.442"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.588"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.338"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06114"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02256"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.424"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.181"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.18%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5839"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9984"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.24"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.07%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.78%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.277"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.382"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +15.42%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5537"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "12.10"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.23%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.928"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07080"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.64"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +21.51%  "
